# Swap the contents of columns C (codeforiati:group-name) and D
# (codeforiati:group-code) for every row of the table, including the
# header row. This matches the IATI codelist re-export which simply
# reorders the group-code / group-name columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
